$p = $ppt.ActivePresentation

# Add a new slide at the end using the "Title and Content" layout (layout index 2
# = "Título e Conteúdo", matching the layout used by most slides in this deck).
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 2)

$s.Shapes.Item(1).TextFrame.TextRange.Text = "Teste"

# Match the deck-wide slide transition convention used by the other slides
# (Transitions > Random Bars, Vertical, Slow).
$s.SlideShowTransition.EntryEffect = 2306   # ppEffectRandomBarsVertical
$s.SlideShowTransition.Speed = 1            # ppTransitionSpeedSlow
